$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maria Yeguez's email in D3 changed from mariavyeguezp@gmail.com
# to mario.calderons@empresa.com - update both the displayed text and
# the underlying mailto hyperlink target, keeping everything else (cell
# style, the other hyperlinks in D2/D4/D7) untouched.
$ws.Range("D3").Value = "mario.calderons@empresa.com"

foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Address($false, $false) -eq "D3") {
        $hl.Address = "mailto:mario.calderons@empresa.com"
    }
}

# Update the selected/active cell to D9 (matches the saved sheet view state)
$ws.Range("D9").Select()
